# Añadidos apellidos en los excel
# Insert a new "Apellido" (surname) column between the existing
# "Candidato"/"Nombre" column and the "Partido" column, and rename the
# first header from "Candidato" to "Nombre".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new column B; this shifts the existing Partido/NIF columns
# from B/C to C/D.
$ws.Columns("B:B").Insert()

# Update header row.
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Apellido"

# Fill in the new surname ("Apellido") column for each candidate row.
$ws.Range("B3").Value = "De Dios"
$ws.Range("B4").Value = "Estevez"
$ws.Range("B5").Value = "Fernandez"
$ws.Range("B6").Value = "Garcia"

# Match the active cell selection recorded in the saved sheet.
$ws.Range("E6").Select()
